$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5541968941688538
$ws.Range("B1").Value = 3.876787900924683
$ws.Range("C1").Value = 5.913495063781738
$ws.Range("D1").Value = 1.444257259368896
$ws.Range("E1").Value = 0.8451204299926758
